# Apply the tracked-requirements status updates described in the commit:
#  - Rows 11, 31, 35, 36 move from "Incomplete"/"In Progress" to "Complete"
#    (green font, shared text "Complete"), and rows 11/30/35/36 gain a
#    completion date in column C (formatted like the other date cells).
#  - The "Picture viewer" requirement text changes from "caption" to "label".
#  - The sheet's remembered selection moves from B32 to C35.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$completeColor = 5287936   # RGB(0,176,80) - same green used by the other "Complete" cells

function Set-Complete($addr) {
    $cell = $ws.Range($addr)
    $cell.Value = "Complete"
    $cell.Font.Color = $completeColor
}

# Row 11: R-FA3 -> Complete, dated 11/2/2011
Set-Complete "B11"
$ws.Range("C11").Value = 40849
$ws.Range("C2").Copy()
$ws.Range("C11").PasteSpecial(-4122)

# Row 30: already Complete, just gains a completion date of 11/3/2011
$ws.Range("C30").Value = 40850
$ws.Range("C2").Copy()
$ws.Range("C30").PasteSpecial(-4122)

# Row 31: R-PV2 -> Complete
Set-Complete "B31"

# Row 35: R-PV7 -> Complete, dated 11/3/2011
Set-Complete "B35"
$ws.Range("C35").Value = 40850
$ws.Range("C2").Copy()
$ws.Range("C35").PasteSpecial(-4122)

# Row 36: R-UP2 -> Complete, dated 11/3/2011
Set-Complete "B36"
$ws.Range("C36").Value = 40850
$ws.Range("C2").Copy()
$ws.Range("C36").PasteSpecial(-4122)

# Wording tweak: "caption" -> "label"
$ws.Range("E35").Value = "Picture viewer will display label displaying the photo’s name"

# Remember the new active selection
$ws.Range("C35").Select()
